# Applies the "refactor pos duration by tcp tool to include component in
# material" edit: changes the pipe material description (polyethylene ->
# cast iron), its diameter, and recomputes every length/duration figure
# that is derived from the new network length (2,2 km -> 3,2 km) further
# down in the document.

function Find-NthRange {
    param($doc, $scopeStart, $scopeEnd, $searchText, $occurrence)
    $count = 0
    $cur = $scopeStart
    while ($cur -lt $scopeEnd) {
        $rng = $doc.Range($cur, $scopeEnd)
        $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
        if (-not $found) { return $null }
        if ($rng.Start -ge $scopeEnd -or $rng.End -gt $scopeEnd) { return $null }
        $count = $count + 1
        if ($count -eq $occurrence) {
            return $rng
        }
        $cur = $rng.End
    }
    return $null
}

function Replace-NthInParagraph {
    param($doc, $paraIndex, $searchText, $replaceText, $occurrence)
    $p = $doc.Paragraphs.Item($paraIndex).Range
    $rng = Find-NthRange $doc $p.Start $p.End $searchText $occurrence
    if ($null -eq $rng) {
        throw "Not found: '$searchText' occurrence $occurrence in paragraph $paraIndex"
    }
    $rng.Text = $replaceText
}

$d = $word.ActiveDocument

# --- Paragraph 2: material/diameter description + total network length ---
Replace-NthInParagraph $d 2 "полиэтиленовых в одну нитку" "чугунных труб" 1
Replace-NthInParagraph $d 2 "труб диаметром " "диаметром " 1
Replace-NthInParagraph $d 2 "до 200" "500" 1
Replace-NthInParagraph $d 2 "полиэтиленовых в одну нитку" "чугунных труб" 1
Replace-NthInParagraph $d 2 "2,2" "3,2" 1

# --- Paragraph 4: "длиной 1 км. составляет 1 мес." -> "длиной 2 км. составляет 4 мес." ---
Replace-NthInParagraph $d 4 "1" "2" 1
Replace-NthInParagraph $d 4 "1" "4" 1

# --- Paragraph 5: "длиной 3 км. составляет 1,5 мес." -> "длиной 4 км. составляет 5 мес." ---
Replace-NthInParagraph $d 5 "3" "4" 1
Replace-NthInParagraph $d 5 "1,5" "5" 1

# --- Paragraph 6 (italic): "длиной 2,2 км:" -> "длиной 3,2 км:" ---
Replace-NthInParagraph $d 6 "2,2" "3,2" 1

# --- Paragraph 8: "(1,5 - 1) / (3 - 1) = 0,2 мес." -> "(5 - 4) / (4 - 2) = 0,5 мес." ---
Replace-NthInParagraph $d 8 "1,5" "5" 1
Replace-NthInParagraph $d 8 "1" "4" 1
Replace-NthInParagraph $d 8 "3" "4" 1
Replace-NthInParagraph $d 8 "1" "2" 1
Replace-NthInParagraph $d 8 "0,2" "0,5" 1

# --- Paragraph 10: "2,2 - 1 = 1,2 км." -> "3,2 - 2 = 1,2 км." (1,2 unchanged) ---
Replace-NthInParagraph $d 10 "2,2" "3,2" 1
Replace-NthInParagraph $d 10 "1" "2" 1

# --- Paragraph 12: "1 + 0,2 ∙ 1,2 = 1,2 мес." -> "4 + 0,5 ∙ 1,2 = 4,6 мес." ---
Replace-NthInParagraph $d 12 "1" "4" 1
Replace-NthInParagraph $d 12 "0,2" "0,5" 1
Replace-NthInParagraph $d 12 "1,2" "4,6" 2

# --- Paragraph 13: "длиной 2,2 км. составляет 1,2 мес." -> "длиной 3,2 км. составляет 4,6 мес." ---
Replace-NthInParagraph $d 13 "2,2" "3,2" 1
Replace-NthInParagraph $d 13 "1,2" "4,6" 1

# --- Paragraph 14: "строительства 1 мес, в т.ч. - 0,1 мес." -> "4,5 мес, в т.ч. - 0,4 мес." ---
Replace-NthInParagraph $d 14 "1" "4,5" 1
Replace-NthInParagraph $d 14 "0,1" "0,4" 1

Write-Host "Done"
